$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text values that look like plain numbers (single decimal point, e.g. "586.29")
# must be forced to Text so Excel does not silently convert them to a numeric value,
# matching the source data which stores every Price/Volume cell as a literal string.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "66.338.19"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "3.317.81"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "586.29"
Set-TextValue $ws.Range("D6") "180.38"
$ws.Range("E6").Value = "  -0.23%  "
Set-TextValue $ws.Range("D7") "0.653"
$ws.Range("E7").Value = "  +5.91%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "3.315.56"
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  -0.62%  "
Set-TextValue $ws.Range("D11") "6.82"
$ws.Range("E11").Value = "  +2.67%  "
Set-TextValue $ws.Range("D12") "0.401"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "3.896.99"
$ws.Range("E13").Value = "  +0.32%  "
Set-TextValue $ws.Range("D14") "0.130"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").Value = "66.377.48"
$ws.Range("E15").Value = "  +0.34%  "
Set-TextValue $ws.Range("D16") "26.52"
$ws.Range("E16").Value = "  -0.34%  "
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "3.285.82"
$ws.Range("E18").Value = "  -0.76%  "
Set-TextValue $ws.Range("D19") "422.98"
$ws.Range("E19").Value = "  -3.23%  "
Set-TextValue $ws.Range("D20") "13.13"
$ws.Range("E20").Value = "  -2.81%  "
Set-TextValue $ws.Range("D21") "5.48"
$ws.Range("E21").Value = "  -3.09%  "
Set-TextValue $ws.Range("D22") "7.35"
$ws.Range("E22").Value = "  -2.53%  "
Set-TextValue $ws.Range("D23") "71.69"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  +0.32%  "
Set-TextValue $ws.Range("D25") "5.67"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "3.470.69"
$ws.Range("E26").Value = "  +0.05%  "
Set-TextValue $ws.Range("D27") "0.514"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("E28").Value = "  +6.59%  "
$ws.Range("E29").Value = "  -0.72%  "
Set-TextValue $ws.Range("D30") "9.10"
$ws.Range("E30").Value = "  +0.73%  "
Set-TextValue $ws.Range("D31") "0.998"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("E32").Value = "  -1.65%  "
Set-TextValue $ws.Range("D33") "22.35"
$ws.Range("E33").Value = "  -1.23%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  -0.67%  "
Set-TextValue $ws.Range("D36") "6.60"
$ws.Range("E36").Value = "  -1.95%  "
Set-TextValue $ws.Range("D37") "1.18"
$ws.Range("E37").Value = "  -2.18%  "
Set-TextValue $ws.Range("D38") "159.84"
$ws.Range("E38").Value = "  +0.04%  "
Set-TextValue $ws.Range("D39") "1.43"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "2.865.90"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D41") "1.80"
$ws.Range("E41").Value = "  +0.63%  "
Set-TextValue $ws.Range("D42") "26.31"
$ws.Range("E42").Value = "  -4.46%  "
Set-TextValue $ws.Range("D43") "4.32"
$ws.Range("E43").Value = "  -2.14%  "
Set-TextValue $ws.Range("D44") "0.758"
$ws.Range("E44").Value = "  -3.64%  "
Set-TextValue $ws.Range("D45") "39.71"
$ws.Range("E45").Value = "  -1.18%  "
Set-TextValue $ws.Range("D46") "0.0659"
$ws.Range("E46").Value = "  -0.55%  "
Set-TextValue $ws.Range("D47") "5.90"
$ws.Range("E47").Value = "  -4.25%  "
Set-TextValue $ws.Range("D48") "2.31"
$ws.Range("E48").Value = "  -1.19%  "
Set-TextValue $ws.Range("D49") "23.10"
$ws.Range("E49").Value = "  -3.80%  "
Set-TextValue $ws.Range("D50") "311.90"
$ws.Range("E50").Value = "  -3.04%  "
Set-TextValue $ws.Range("D51") "0.0272"
$ws.Range("E51").Value = "  +0.51%  "
